# Apply the updated cryptocurrency price/volume snapshot to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text, $numericLooking) {
    $r = $ws.Range($cellRef)
    if ($numericLooking) {
        # Leading apostrophe forces Excel to keep a numeric-looking value as text,
        # matching the source data (prices/IDs stored as literal strings).
        $r.Value = "'" + $text
        $r.Style = "Normal"
    } else {
        $r.Value = $text
    }
}

$updates = @(
    @{Cell = "D2"; Text = '59.411.50'; Numeric = $false},
    @{Cell = "E2"; Text = '  -1.72%  '; Numeric = $false},
    @{Cell = "D3"; Text = '2.580.28'; Numeric = $false},
    @{Cell = "E3"; Text = '  -2.21%  '; Numeric = $false},
    @{Cell = "E4"; Text = '  -0.29%  '; Numeric = $false},
    @{Cell = "D5"; Text = '556.57'; Numeric = $true},
    @{Cell = "E5"; Text = '  -1.88%  '; Numeric = $false},
    @{Cell = "D7"; Text = '0.999'; Numeric = $true},
    @{Cell = "E7"; Text = '  +0.13%  '; Numeric = $false},
    @{Cell = "E8"; Text = '  -2.41%  '; Numeric = $false},
    @{Cell = "D9"; Text = '2.587.08'; Numeric = $false},
    @{Cell = "E9"; Text = '  -2.86%  '; Numeric = $false},
    @{Cell = "D10"; Text = '6.67'; Numeric = $true},
    @{Cell = "E10"; Text = '  -2.37%  '; Numeric = $false},
    @{Cell = "E11"; Text = '  -0.62%  '; Numeric = $false},
    @{Cell = "E12"; Text = '  +12.77%  '; Numeric = $false},
    @{Cell = "D13"; Text = '0.353'; Numeric = $true},
    @{Cell = "E13"; Text = '  +2.86%  '; Numeric = $false},
    @{Cell = "D14"; Text = '3.036.17'; Numeric = $false},
    @{Cell = "E14"; Text = '  -2.57%  '; Numeric = $false},
    @{Cell = "D15"; Text = '59.400.28'; Numeric = $false},
    @{Cell = "E15"; Text = '  -1.73%  '; Numeric = $false},
    @{Cell = "D16"; Text = '23.06'; Numeric = $true},
    @{Cell = "E16"; Text = '  +5.13%  '; Numeric = $false},
    @{Cell = "E17"; Text = '  -0.56%  '; Numeric = $false},
    @{Cell = "D18"; Text = '2.586.92'; Numeric = $false},
    @{Cell = "E18"; Text = '  -2.34%  '; Numeric = $false},
    @{Cell = "E19"; Text = '  +0.21%  '; Numeric = $false},
    @{Cell = "D20"; Text = '337.30'; Numeric = $true},
    @{Cell = "E20"; Text = '  -1.60%  '; Numeric = $false},
    @{Cell = "D21"; Text = '10.35'; Numeric = $true},
    @{Cell = "E21"; Text = '  -0.75%  '; Numeric = $false},
    @{Cell = "D22"; Text = '6.43'; Numeric = $true},
    @{Cell = "E22"; Text = '  +1.08%  '; Numeric = $false},
    @{Cell = "E23"; Text = '  -0.01%  '; Numeric = $false},
    @{Cell = "D24"; Text = '63.20'; Numeric = $true},
    @{Cell = "E24"; Text = '  -4.41%  '; Numeric = $false},
    @{Cell = "E25"; Text = '  +6.72%  '; Numeric = $false},
    @{Cell = "E26"; Text = '  +0.26%  '; Numeric = $false},
    @{Cell = "E27"; Text = '  -2.33%  '; Numeric = $false},
    @{Cell = "D28"; Text = '7.44'; Numeric = $true},
    @{Cell = "E28"; Text = '  +0.56%  '; Numeric = $false},
    @{Cell = "E29"; Text = '  -3.36%  '; Numeric = $false},
    @{Cell = "E30"; Text = '  -0.03%  '; Numeric = $false},
    @{Cell = "D31"; Text = '6.18'; Numeric = $true},
    @{Cell = "E31"; Text = '  -0.79%  '; Numeric = $false},
    @{Cell = "D32"; Text = '1.66'; Numeric = $true},
    @{Cell = "E32"; Text = '  -2.76%  '; Numeric = $false},
    @{Cell = "D33"; Text = '157.87'; Numeric = $true},
    @{Cell = "E33"; Text = '  -0.98%  '; Numeric = $false},
    @{Cell = "D34"; Text = '19.09'; Numeric = $true},
    @{Cell = "E34"; Text = '  -0.69%  '; Numeric = $false},
    @{Cell = "D35"; Text = '4.05'; Numeric = $true},
    @{Cell = "E35"; Text = '  -0.82%  '; Numeric = $false},
    @{Cell = "B36"; Text = 'Fetch.AI'; Numeric = $false},
    @{Cell = "C36"; Text = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; Numeric = $false},
    @{Cell = "D36"; Text = '0.912'; Numeric = $true},
    @{Cell = "E36"; Text = '  +0.82%  '; Numeric = $false},
    @{Cell = "B37"; Text = 'ImmutableX'; Numeric = $false},
    @{Cell = "C37"; Text = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Numeric = $false},
    @{Cell = "D37"; Text = '1.16'; Numeric = $true},
    @{Cell = "E37"; Text = '  -0.20%  '; Numeric = $false},
    @{Cell = "D38"; Text = '37.42'; Numeric = $true},
    @{Cell = "E38"; Text = '  -0.16%  '; Numeric = $false},
    @{Cell = "D39"; Text = '0.859'; Numeric = $true},
    @{Cell = "E39"; Text = '  -5.49%  '; Numeric = $false},
    @{Cell = "E40"; Text = '  -1.91%  '; Numeric = $false},
    @{Cell = "D41"; Text = '3.67'; Numeric = $true},
    @{Cell = "E41"; Text = '  +0.54%  '; Numeric = $false},
    @{Cell = "D42"; Text = '292.18'; Numeric = $true},
    @{Cell = "E42"; Text = '  -3.85%  '; Numeric = $false},
    @{Cell = "D43"; Text = '136.50'; Numeric = $true},
    @{Cell = "E43"; Text = '  +5.84%  '; Numeric = $false},
    @{Cell = "E44"; Text = '  +0.42%  '; Numeric = $false},
    @{Cell = "D45"; Text = '0.0976'; Numeric = $true},
    @{Cell = "E45"; Text = '  -1.09%  '; Numeric = $false},
    @{Cell = "D46"; Text = '0.592'; Numeric = $true},
    @{Cell = "E46"; Text = '  -2.04%  '; Numeric = $false},
    @{Cell = "D47"; Text = '10.67'; Numeric = $true},
    @{Cell = "E48"; Text = '  -2.27%  '; Numeric = $false},
    @{Cell = "D49"; Text = '0.0234'; Numeric = $true},
    @{Cell = "E49"; Text = '  -1.27%  '; Numeric = $false},
    @{Cell = "D50"; Text = '18.78'; Numeric = $true},
    @{Cell = "E50"; Text = '  -0.01%  '; Numeric = $false},
    @{Cell = "D51"; Text = '1.957.07'; Numeric = $false},
    @{Cell = "E51"; Text = '  -0.23%  '; Numeric = $false}
)

foreach ($u in $updates) {
    Set-CellText $u.Cell $u.Text $u.Numeric
}
